$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data values
$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = " RESIDENT-1"
$ws.Range("C2").Value = "Resident Virtual Machine"
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "RESIDENT-REG"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "Resident Virtual Machine"
$ws.Range("I2").Value = $true

# Custom boolean display format for is_active column
$ws.Range("I2").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# Default column width for the sheet
$ws.StandardWidth = 11.53515625

# Column widths (in character units, matching target dimension)
$ws.Columns.Item(1).ColumnWidth = 10.31
$ws.Columns.Item(2).ColumnWidth = 16.2
$ws.Columns.Item(3).ColumnWidth = 30.81
$ws.Columns.Item(4).ColumnWidth = 7.01
$ws.Columns.Item(5).ColumnWidth = 7.22
$ws.Columns.Item(6).ColumnWidth = 14.4
$ws.Columns.Item(7).ColumnWidth = 17.7
$ws.Columns.Item(8).ColumnWidth = 24.9
$ws.Columns.Item(9).ColumnWidth = 9.02

# Selection to match target
$ws.Range("B2").Select()
